$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.507882
$ws.Range("H2").Value = 16.523646
$ws.Range("I2").Value = 0.03518866199235487
$ws.Range("J2").Value = 0.03518866199235487
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 20.545366
$ws.Range("N2").Value = 61.636098
$ws.Range("O2").Value = 0.8965950288338865
$ws.Range("P2").Value = 0.8965950288338865
$ws.Range("Q2").Value = 113.161451574812
$ws.Range("R2").Value = 1018.453064173308
$ws.Range("S2").Value = 0.0315499794136613
$ws.Range("T2").Value = 0.0315499794136613
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.507882
$ws.Range("H3").Value = 16.523646
$ws.Range("I3").Value = 0.03518866199235487
$ws.Range("J3").Value = 0.03518866199235487
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.154739666666666
$ws.Range("N3").Value = 6.464219
$ws.Range("O3").Value = 0.09403234157836461
$ws.Range("P3").Value = 0.09403234157836463
$ws.Range("Q3").Value = 11.86805182471933
$ws.Range("R3").Value = 106.812466422474
$ws.Range("S3").Value = 0.003308872284150729
$ws.Range("T3").Value = 0.003308872284150729
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.507882
$ws.Range("H4").Value = 16.523646
$ws.Range("I4").Value = 0.03518866199235487
$ws.Range("J4").Value = 0.03518866199235487
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2147726666666667
$ws.Range("N4").Value = 0.6443180000000001
$ws.Range("O4").Value = 0.009372629587748921
$ws.Range("P4").Value = 0.009372629587748921
$ws.Range("Q4").Value = 1.182942504825333
$ws.Range("R4").Value = 10.646482543428
$ws.Range("S4").Value = 0.0003298102945428411
$ws.Range("T4").Value = 0.0003298102945428411
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 149.9875183333334
$ws.Range("H5").Value = 449.9625550000001
$ws.Range("I5").Value = 0.9582376829612175
$ws.Range("J5").Value = 0.9582376829612176
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 20.545366
$ws.Range("N5").Value = 61.636098
$ws.Range("O5").Value = 0.8965950288338865
$ws.Range("P5").Value = 0.8965950288338865
$ws.Range("Q5").Value = 3081.548459590044
$ws.Range("R5").Value = 27733.9361363104
$ws.Range("S5").Value = 0.8591511429843294
$ws.Range("T5").Value = 0.8591511429843295
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 149.9875183333334
$ws.Range("H6").Value = 449.9625550000001
$ws.Range("I6").Value = 0.9582376829612175
$ws.Range("J6").Value = 0.9582376829612176
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.154739666666666
$ws.Range("N6").Value = 6.464219
$ws.Range("O6").Value = 0.09403234157836461
$ws.Range("P6").Value = 0.09403234157836463
$ws.Range("Q6").Value = 323.1840552577273
$ws.Range("R6").Value = 2908.656497319545
$ws.Range("S6").Value = 0.09010533311746986
$ws.Range("T6").Value = 0.09010533311746988
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 149.9875183333334
$ws.Range("H7").Value = 449.9625550000001
$ws.Range("I7").Value = 0.9582376829612175
$ws.Range("J7").Value = 0.9582376829612176
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.2147726666666667
$ws.Range("N7").Value = 0.6443180000000001
$ws.Range("O7").Value = 0.009372629587748921
$ws.Range("P7").Value = 0.009372629587748921
$ws.Range("Q7").Value = 32.21321927916556
$ws.Range("R7").Value = 289.9189735124901
$ws.Range("S7").Value = 0.008981206859418276
$ws.Range("T7").Value = 0.008981206859418278
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.028937
$ws.Range("H8").Value = 3.086811
$ws.Range("I8").Value = 0.006573655046427582
$ws.Range("J8").Value = 0.006573655046427582
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 20.545366
$ws.Range("N8").Value = 61.636098
$ws.Range("O8").Value = 0.8965950288338865
$ws.Range("P8").Value = 0.8965950288338865
$ws.Range("Q8").Value = 21.139887255942
$ws.Range("R8").Value = 190.258985303478
$ws.Range("S8").Value = 0.005893906435895761
$ws.Range("T8").Value = 0.005893906435895761
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.028937
$ws.Range("H9").Value = 3.086811
$ws.Range("I9").Value = 0.006573655046427582
$ws.Range("J9").Value = 0.006573655046427582
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.154739666666666
$ws.Range("N9").Value = 6.464219
$ws.Range("O9").Value = 0.09403234157836461
$ws.Range("P9").Value = 0.09403234157836463
$ws.Range("Q9").Value = 2.217091368401
$ws.Range("R9").Value = 19.953822315609
$ws.Range("S9").Value = 0.0006181361767440187
$ws.Range("T9").Value = 0.0006181361767440188
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.028937
$ws.Range("H10").Value = 3.086811
$ws.Range("I10").Value = 0.006573655046427582
$ws.Range("J10").Value = 0.006573655046427582
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.2147726666666667
$ws.Range("N10").Value = 0.6443180000000001
$ws.Range("O10").Value = 0.009372629587748921
$ws.Range("P10").Value = 0.009372629587748921
$ws.Range("Q10").Value = 0.220987543322
$ws.Range("R10").Value = 1.988887889898
$ws.Range("S10").Value = 0.00006161243378780217
$ws.Range("T10").Value = 0.00006161243378780217

Write-Output "Updated cells"